# Update countries & provincias Spain
#
# 1) Swap the displayed country names for rows 206/207 (Santa Lucia now
#    comes before Timor Oriental in the shared-strings list) while the
#    rest of each row's data stays put.
# 2) Bump the "Datos actualizados..." timestamp string.
# 3) Refresh the daily case-count figures for a handful of countries
#    (Estados Unidos, Brasil, Alemania, Canada, Curazao).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Santa Lucia / Timor Oriental swap -------------------------------
$ws.Range("A206").Value = "Santa Lucia"
$ws.Range("A207").Value = "Timor Oriental"

# --- 2) Updated timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 26 de Septiembre de 2020 a las 23:33"

# --- 3) Refreshed country statistics --------------------------------------
# Estados Unidos (row 4)
$ws.Range("B4").Value = 7281459
$ws.Range("C4").Value = 37275
$ws.Range("D4").Value = 4512301
$ws.Range("E4").Value = 2560066
$ws.Range("G4").Value = 652
$ws.Range("H4").Value = 209092

# Brasil (row 6)
$ws.Range("B6").Value = 4717991
$ws.Range("C6").Value = 25412
$ws.Range("E6").Value = 535636
$ws.Range("G6").Value = 697
$ws.Range("H6").Value = 141406

# Alemania (row 25)
$ws.Range("B25").Value = 285025
$ws.Range("C25").Value = 1319
$ws.Range("E25").Value = 25993

# Canada (row 29)
$ws.Range("B29").Value = 151671
$ws.Range("C29").Value = 1215
$ws.Range("D29").Value = 130328
$ws.Range("E29").Value = 12080

# Curazao (row 185)
$ws.Range("B185").Value = 337
$ws.Range("C185").Value = 8
$ws.Range("D185").Value = 134
$ws.Range("E185").Value = 202
